# "new tibble update adjustments"
# A new row-2 "-" placeholder is added in columns A and B (a fresh shared
# string), and the sheet's active selection moves from A4:XFD4 to B3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "-"
$ws.Range("B2").Value = "-"

$ws.Range("B3").Select()
